$d = $word.ActiveDocument

$p = $d.Paragraphs(1)

# Replace paragraph text (minus the trailing paragraph mark) so that the
# two runs ("**ID__AFFARS_pgi_5334_topic_3__ID**" + trailing space) collapse
# into a single run carrying the new identifier and no trailing space.
$r = $p.Range
[void]$r.MoveEnd(1, -1)
$r.Text = "**ID__AFFARS_SMC_PGI_5334_2__ID**"

# Add a paragraph border (padding-only, 5 twips on every side).
$p.Borders.DistanceFromTop = 5
$p.Borders.DistanceFromLeft = 5
$p.Borders.DistanceFromBottom = 5
$p.Borders.DistanceFromRight = 5

# Widen the left indent (LeftIndent is in points; 225 twips = 11.25 pt).
$p.LeftIndent = 11.25
